# Update Name of Algo
# Apply updated imputed values to column D (algorithm result column) for
# the specified rows, matching the canonical OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value  = -7.575
$ws.Range("D10").Value = -7.999
$ws.Range("D12").Value = -7.937
$ws.Range("D18").Value = -8.038999999999998
$ws.Range("D25").Value = -8.183
